# Swap columns B and C (header text and all data values) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2

    $bCell.Value2 = $cVal
    $cCell.Value2 = $bVal
}
